$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 45: "get conversation history ( POST )" Voice Bot API -------

# Column A: API name
$ws.Range("A45").Value = "get conversation history ( POST )"

# Column B: URL, rendered as a hyperlink (matches the style used by the
# other API rows further up the sheet, e.g. B39 / B42)
$ws.Range("B45").Value = "http://1msg.1point1.in:3001/api/auth/j-v1/get/conversation/history/"
[void]$ws.Hyperlinks.Add($ws.Range("B45"), "http://1msg.1point1.in:3001/api/auth/j-v1/get/conversation/history/")
$ws.Range("B45").Style = "Hyperlink"

# Column C: JSON payload (wrapped, multi-line)
$payload = "{`n    ""user_id"" : 11,`n    ""agent_id"": 72 ,`n    ""from_date"" : ""2025-05-18"" ,`n    ""to_date"" : ""2025-05-18""`n}"
$ws.Range("C45").Value = $payload
$ws.Range("C45").WrapText = $true

# Column D: curl command reproducing the same request
$curl = "curl --location 'http://1msg.1point1.in:3001/api/auth/j-v1/get/conversation/history/' \`n--header 'Content-Type: application/json' \`n--data '{`n    ""user_id"" : 11,`n    ""agent_id"": 72 ,`n    ""from_date"" : ""2025-05-18"" ,`n    ""to_date"" : ""2025-05-18""`n}'"
$ws.Range("D45").Value = $curl
$ws.Range("D45").WrapText = $true

# Row height matches the tall, wrapped, multi-line content
$ws.Rows.Item(45).RowHeight = 124.2

# Leave the active selection on D45, matching the authored selection after
# the edit (the view also scrolls so the new row is visible)
[void]$ws.Range("A43").Select()
[void]$ws.Range("D45").Select()
$excel.ActiveWindow.ScrollRow = 43
